$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 576, shifting existing rows 576:635 down to 577:636
$ws.Rows.Item(576).Insert()

# Populate the new row 576 with the new weekly data entry
$ws.Cells.Item(576, 1).Value = 10
$ws.Cells.Item(576, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(576, 3).Value = "La Araucanía"
$ws.Cells.Item(576, 4).Value = 45166
$ws.Cells.Item(576, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(576, 5).Value = 9
$ws.Cells.Item(576, 6).Value = 100112040
$ws.Cells.Item(576, 7).Value = "Cilantro"
$ws.Cells.Item(576, 8).Value = "Sin especificar"
$ws.Cells.Item(576, 9).Value = "Primera"
$ws.Cells.Item(576, 10).Value = 60
$ws.Cells.Item(576, 11).Value = 5000
$ws.Cells.Item(576, 12).Value = 5000
$ws.Cells.Item(576, 13).Value = 5000
$ws.Cells.Item(576, 14).Value = "`$/docena de atados (2 kilos)"
$ws.Cells.Item(576, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(576, 16).Value = 2500
$ws.Cells.Item(576, 17).Value = 2
$ws.Cells.Item(576, 18).Value = "Hortaliza"
